# -----------------------------------------------------------------------
# "cryptos" price-table refresh (Sheet1, columns B:E, rows 2-51).
# Almost every Price (D) / Volume(1h) (E) cell gets a new scraped value;
# two coin pairs also swapped row position while keeping the rank index in
# column A untouched: Uniswap/BitcoinCash (rows 19-20) and
# ImmutableX/Maker (rows 37-38).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val, [bool]$numericLooking) {
    $cell = $ws.Range($addr)
    if ($numericLooking) {
        # Plain decimal-looking price strings (e.g. "214.08", "0.05160") would
        # otherwise get auto-coerced by Excel into a Double on assignment,
        # silently dropping significant trailing zeros / introducing float
        # noise. Force Text first, write the literal, then drop the now-
        # redundant explicit format so the cell ends up as styleless as the
        # rest of the sheet.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.ClearFormats()
    } else {
        $cell.Value = $val
    }
}

Set-TextCell "D2" '25.912.46' $false
Set-TextCell "E2" '  -0.20%  ' $false
Set-TextCell "D3" '1.635.74' $false
Set-TextCell "E3" '  -0.42%  ' $false
Set-TextCell "E4" '  -0.43%  ' $false
Set-TextCell "D5" '214.08' $true
Set-TextCell "E5" '  -0.74%  ' $false
Set-TextCell "D6" '0.5054' $true
Set-TextCell "E6" '  -0.59%  ' $false
Set-TextCell "E8" '  +0.22%  ' $false
Set-TextCell "D9" '0.06356' $true
Set-TextCell "E10" '  +0.79%  ' $false
Set-TextCell "D11" '0.07730' $true
Set-TextCell "E11" '  -0.75%  ' $false
Set-TextCell "D12" '4.277' $true
Set-TextCell "E12" '  -0.52%  ' $false
Set-TextCell "D13" '1.635.62' $false
Set-TextCell "E13" '  -1.03%  ' $false
Set-TextCell "D14" '0.5431' $true
Set-TextCell "E14" '  -0.84%  ' $false
Set-TextCell "D15" '0.0₅7727' $false
Set-TextCell "E15" '  -1.70%  ' $false
Set-TextCell "D16" '64.01' $true
Set-TextCell "E16" '  -0.60%  ' $false
Set-TextCell "D17" '25.925.04' $false
Set-TextCell "E17" '  -0.37%  ' $false
Set-TextCell "E18" '  -0.24%  ' $false
Set-TextCell "B19" 'Uniswap' $false
Set-TextCell "C19" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' $false
Set-TextCell "D19" '4.427' $true
Set-TextCell "E19" '  -0.59%  ' $false
Set-TextCell "B20" 'BitcoinCash' $false
Set-TextCell "C20" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' $false
Set-TextCell "D20" '194.86' $true
Set-TextCell "E20" '  -1.73%  ' $false
Set-TextCell "D21" '9.903' $true
Set-TextCell "E21" '  -0.77%  ' $false
Set-TextCell "D22" '6.098' $true
Set-TextCell "E22" '  +0.50%  ' $false
Set-TextCell "E23" '  -0.34%  ' $false
Set-TextCell "D24" '1.890' $true
Set-TextCell "E24" '  +0.82%  ' $false
Set-TextCell "D25" '143.06' $true
Set-TextCell "E25" '  +1.21%  ' $false
Set-TextCell "D26" '0.1239' $true
Set-TextCell "E26" '  +7.15%  ' $false
Set-TextCell "D27" '6.812' $true
Set-TextCell "E27" '  -1.03%  ' $false
Set-TextCell "E28" '  -0.89%  ' $false
Set-TextCell "E29" '  -0.49%  ' $false
Set-TextCell "D30" '0.04860' $true
Set-TextCell "E30" '  -3.49%  ' $false
Set-TextCell "D31" '3.242' $true
Set-TextCell "E31" '  -0.73%  ' $false
Set-TextCell "E32" '  -0.08%  ' $false
Set-TextCell "D33" '1.550' $true
Set-TextCell "E33" '  +0.33%  ' $false
Set-TextCell "D34" '2.373' $true
Set-TextCell "E34" '  +0.37%  ' $false
Set-TextCell "D35" '0.9102' $true
Set-TextCell "E35" '  +1.10%  ' $false
Set-TextCell "D36" '2.571' $true
Set-TextCell "E36" '  -0.60%  ' $false
Set-TextCell "B37" 'ImmutableX' $false
Set-TextCell "C37" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' $false
Set-TextCell "D37" '0.5495' $true
Set-TextCell "E37" '  -0.24%  ' $false
Set-TextCell "B38" 'Maker' $false
Set-TextCell "C38" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' $false
Set-TextCell "D38" '1.123.88' $false
Set-TextCell "E38" '  -1.13%  ' $false
Set-TextCell "E39" '  -0.20%  ' $false
Set-TextCell "E40" '  -0.36%  ' $false
Set-TextCell "D41" '5.579' $true
Set-TextCell "E41" '  -0.83%  ' $false
Set-TextCell "D42" '0.8022' $true
Set-TextCell "E42" '  -2.01%  ' $false
Set-TextCell "E43" '  -8.64%  ' $false
Set-TextCell "D44" '98.42' $true
Set-TextCell "E44" '  -1.88%  ' $false
Set-TextCell "D45" '1.770.94' $false
Set-TextCell "E45" '  -0.47%  ' $false
Set-TextCell "D46" '0.4474' $true
Set-TextCell "E46" '  -1.35%  ' $false
Set-TextCell "D47" '1.002' $true
Set-TextCell "E47" '  -0.05%  ' $false
Set-TextCell "E48" '  -0.02%  ' $false
Set-TextCell "D49" '0.05160' $true
Set-TextCell "D50" '7.525' $true
Set-TextCell "E50" '  +2.03%  ' $false
Set-TextCell "D51" '0.9990' $true
Set-TextCell "E51" '  -0.72%  ' $false
